# Updated cryptos list on Sat Jun 22 18:29:49 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns stay as plain text so numeric-looking
# strings (e.g. "7.25", "1.00") are not silently converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# --- Row-level value updates (Price column D, Volume(1h) column E) ---

# Row 2: Bitcoin
$ws.Range("D2").Value = "64.262.85"
$ws.Range("E2").Value = "  +0.67%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "3.497.29"
$ws.Range("E3").Value = "  +0.44%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.09%  "

# Row 5: BNB
$ws.Range("D5").Value = "586.15"
$ws.Range("E5").Value = "  +0.61%  "

# Row 6: Solana
$ws.Range("D6").Value = "134.03"
$ws.Range("E6").Value = "  +3.37%  "

# Row 7: USDC
$ws.Range("E7").Value = "  -0.06%  "

# Row 8: XRP
$ws.Range("E8").Value = "  +1.23%  "

# Row 9: Dogecoin
$ws.Range("E9").Value = "  +1.17%  "

# Row 10: Toncoin
$ws.Range("D10").Value = "7.25"
$ws.Range("E10").Value = "  +2.58%  "

# Row 11: Cardano
$ws.Range("E11").Value = "  +3.02%  "

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "4.093.37"
$ws.Range("E12").Value = "  +0.20%  "

# Row 13: ShibaInu
$ws.Range("D13").Value = "0.0000181"
$ws.Range("E13").Value = "  +3.34%  "

# Row 14: TRON
$ws.Range("E14").Value = "  +1.24%  "

# Row 15: WrappedEther
$ws.Range("D15").Value = "3.498.33"
$ws.Range("E15").Value = "  +0.00%  "

# Row 16: Avalanche
$ws.Range("D16").Value = "25.83"
$ws.Range("E16").Value = "  -4.67%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "64.263.32"
$ws.Range("E17").Value = "  +0.88%  "

# Row 18: Uniswap
$ws.Range("E18").Value = "  +1.78%  "

# Row 19: Polkadot
$ws.Range("D19").Value = "5.74"
$ws.Range("E19").Value = "  +2.79%  "

# Row 20: Chainlink
$ws.Range("D20").Value = "13.61"
$ws.Range("E20").Value = "  -2.92%  "

# Row 21: BitcoinCash
$ws.Range("D21").Value = "393.55"
$ws.Range("E21").Value = "  +3.64%  "

# Row 22: Polygon
$ws.Range("D22").Value = "0.570"
$ws.Range("E22").Value = "  +0.32%  "

# Row 23: WrappedeETH
$ws.Range("D23").Value = "3.637.56"
$ws.Range("E23").Value = "  +0.40%  "

# Row 24: Litecoin
$ws.Range("D24").Value = "74.33"
$ws.Range("E24").Value = "  +1.80%  "

# Row 25: Dai
$ws.Range("E25").Value = "  +0.01%  "

# Row 26: PEPE
$ws.Range("E26").Value = "  +0.29%  "

# Row 27: Binance-PegBSC-USD
$ws.Range("E27").Value = "  +0.10%  "

# Row 28: RenderToken
$ws.Range("D28").Value = "7.39"
$ws.Range("E28").Value = "  +0.26%  "

# Row 29: Fetch.AI
$ws.Range("E29").Value = "  -2.90%  "

# Row 30 <-> Row 31: coins swapped places (InternetComputer(DFINITY) and PancakeSwap)
# Row 30 becomes PancakeSwap
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "2.24"
$ws.Range("E30").Value = "  +1.33%  "

# Row 31 becomes InternetComputer(DFINITY)
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "8.26"
$ws.Range("E31").Value = "  +1.24%  "

# Row 32: RenzoRestakedETH
$ws.Range("D32").Value = "3.517.02"
$ws.Range("E32").Value = "  +0.74%  "

# Row 33: Kaspa
$ws.Range("E33").Value = "  +5.19%  "

# Row 34: USDe
$ws.Range("E34").Value = "  +0.05%  "

# Row 35: EthereumClassic
$ws.Range("D35").Value = "23.42"
$ws.Range("E35").Value = "  +0.66%  "

# Row 36: NEARProtocol
$ws.Range("D36").Value = "5.15"
$ws.Range("E36").Value = "  -1.75%  "

# Row 37: Aptos
$ws.Range("D37").Value = "6.89"
$ws.Range("E37").Value = "  +1.09%  "

# Row 38: ImmutableX
$ws.Range("E38").Value = "  +1.05%  "

# Row 39: Monero
$ws.Range("D39").Value = "165.57"
$ws.Range("E39").Value = "  +4.04%  "

# Row 40: Hedera
$ws.Range("D40").Value = "0.0781"
$ws.Range("E40").Value = "  -0.61%  "

# Row 41: Mantle
$ws.Range("E41").Value = "  -0.21%  "

# Row 42: FirstDigitalUSD
$ws.Range("E42").Value = "  -0.21%  "

# Row 43: EnergySwap
$ws.Range("D43").Value = "25.33"
$ws.Range("E43").Value = "  -2.08%  "

# Row 44: Filecoin
$ws.Range("D44").Value = "4.39"
$ws.Range("E44").Value = "  +1.22%  "

# Row 45: ONDO
$ws.Range("E45").Value = "  -1.32%  "

# Row 46: Stacks
$ws.Range("E46").Value = "  +3.50%  "

# Row 47: Maker
$ws.Range("D47").Value = "2.460.91"
$ws.Range("E47").Value = "  +2.09%  "

# Row 48: Cosmos
$ws.Range("E48").Value = "  -0.13%  "

# Row 49: SuiNetwork
$ws.Range("D49").Value = "0.893"
$ws.Range("E49").Value = "  +0.83%  "

# Row 50: VeChain
$ws.Range("E50").Value = "  +0.08%  "

# Row 51: TheGraph -> InjectiveProtocol
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "21.13"
$ws.Range("E51").Value = "  +0.47%  "

# Restore the default (Normal) style on the touched columns so that no
# stray explicit cell styles are introduced by the temporary text format.
$ws.Range("D2:E51").Style = "Normal"
